$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 297, shifting existing rows 297:325 down to 298:326
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row 297 with the new record
$ws.Cells.Item(297, 1).Value = 3
$ws.Cells.Item(297, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(297, 3).Value = "Coquimbo"
$ws.Cells.Item(297, 4).Value = 44769
$ws.Cells.Item(297, 5).Value = 5
$ws.Cells.Item(297, 6).Value = 100112001
$ws.Cells.Item(297, 7).Value = "Berenjena"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 115
$ws.Cells.Item(297, 11).Value = 8000
$ws.Cells.Item(297, 12).Value = 9000
$ws.Cells.Item(297, 13).Value = 8478
$ws.Cells.Item(297, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(297, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(297, 16).Value = 141
$ws.Cells.Item(297, 17).Value = 60
$ws.Cells.Item(297, 18).Value = "Hortaliza"
